$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 688.73334
$ws.Range("I80").Value = 405.33334
$ws.Range("J80").Value = 877.6667
$ws.Range("K80").Value = 1216.00002
$ws.Range("L80").Value = 2633.0001
$ws.Range("M80").Value = -218.0000199999999
$ws.Range("N80").Value = -4629.0001

$ws.Range("H83").Value = 688.73334
$ws.Range("I83").Value = 405.33334
$ws.Range("J83").Value = 877.6667
$ws.Range("K83").Value = 3648.00006
$ws.Range("L83").Value = 7899.0003
$ws.Range("M83").Value = 1343.99994
$ws.Range("N83").Value = -17883.0003

$ws.Range("H138").Value = 1904.79
$ws.Range("I138").Value = 974.8889
$ws.Range("J138").Value = 1996.7583
$ws.Range("K138").Value = 2924.6667
$ws.Range("L138").Value = 5990.2749
$ws.Range("M138").Value = 2215.3333
$ws.Range("N138").Value = -16270.2749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 133.33333
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -424

$ws.Range("H25").Value = 800
$ws.Range("I25").Value = 800
$ws.Range("K25").Value = 800
$ws.Range("M25").Value = -398

$ws.Range("H35").Value = 2548.4285
$ws.Range("I35").Value = 1567.4
$ws.Range("J35").Value = 5001
$ws.Range("K35").Value = 1567.4
$ws.Range("L35").Value = 5001
$ws.Range("M35").Value = -1161.4
$ws.Range("N35").Value = -5813

$ws.Range("H74").Value = 36001200
$ws.Range("I74").Value = 83667250
$ws.Range("J74").Value = 251662.5
$ws.Range("K74").Value = 83667250
$ws.Range("L74").Value = 251662.5
$ws.Range("M74").Value = -83666376
$ws.Range("N74").Value = -253410.5

$ws.Range("H75").Value = 39123.8
$ws.Range("I75").Value = 35000
$ws.Range("J75").Value = 40154.75
$ws.Range("K75").Value = 35000
$ws.Range("L75").Value = 40154.75
$ws.Range("M75").Value = -34126
$ws.Range("N75").Value = -41902.75

$ws.Range("H77").Value = 36001200
$ws.Range("I77").Value = 83667250
$ws.Range("J77").Value = 251662.5
$ws.Range("K77").Value = 418336250
$ws.Range("L77").Value = 1258312.5
$ws.Range("M77").Value = -418331882
$ws.Range("N77").Value = -1267048.5

$ws.Range("H78").Value = 39123.8
$ws.Range("I78").Value = 35000
$ws.Range("J78").Value = 40154.75
$ws.Range("K78").Value = 105000
$ws.Range("L78").Value = 120464.25
$ws.Range("M78").Value = -100632
$ws.Range("N78").Value = -129200.25

$ws.Range("H80").Value = 44433.4
$ws.Range("I80").Value = 35000
$ws.Range("J80").Value = 46791.75
$ws.Range("K80").Value = 35000
$ws.Range("L80").Value = 46791.75
$ws.Range("M80").Value = -34002
$ws.Range("N80").Value = -48787.75

$ws.Range("H81").Value = 40135.75
$ws.Range("J81").Value = 40135.75
$ws.Range("L81").Value = 40135.75
$ws.Range("N81").Value = -42131.75

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H83").Value = 44433.4
$ws.Range("I83").Value = 35000
$ws.Range("J83").Value = 46791.75
$ws.Range("K83").Value = 105000
$ws.Range("L83").Value = 140375.25
$ws.Range("M83").Value = -100008
$ws.Range("N83").Value = -150359.25

$ws.Range("H84").Value = 40135.75
$ws.Range("J84").Value = 40135.75
$ws.Range("L84").Value = 120407.25
$ws.Range("N84").Value = -130391.25

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H86").Value = 39250.4
$ws.Range("I86").Value = 35000
$ws.Range("J86").Value = 40313
$ws.Range("K86").Value = 35000
$ws.Range("L86").Value = 40313
$ws.Range("M86").Value = -33814
$ws.Range("N86").Value = -42685

$ws.Range("H87").Value = 29175.5
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 30195
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 30195
$ws.Range("M87").Value = -18752
$ws.Range("N87").Value = -32691

$ws.Range("H89").Value = 39250.4
$ws.Range("I89").Value = 35000
$ws.Range("J89").Value = 40313
$ws.Range("K89").Value = 105000
$ws.Range("L89").Value = 120939
$ws.Range("M89").Value = -99072
$ws.Range("N89").Value = -132795

$ws.Range("H90").Value = 29175.5
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 30195
$ws.Range("K90").Value = 60000
$ws.Range("L90").Value = 90585
$ws.Range("M90").Value = -53760
$ws.Range("N90").Value = -103065

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H110").Value = 358104.6
$ws.Range("I110").Value = 556677.4
$ws.Range("J110").Value = 673.6
$ws.Range("K110").Value = 556677.4
$ws.Range("L110").Value = 673.6
$ws.Range("M110").Value = -554632.4
$ws.Range("N110").Value = -4763.6

$ws.Range("H132").Value = 107644.63
$ws.Range("I132").Value = 73381.14
$ws.Range("K132").Value = 220143.42
$ws.Range("M132").Value = -217613.42

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 133.33333
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -430

$ws.Range("H24").Value = 2700
$ws.Range("I24").Value = 400
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 400
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -165
$ws.Range("N24").Value = -5470

$ws.Range("H75").Value = 11457
$ws.Range("I75").Value = 2914
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 2914
$ws.Range("L75").Value = 20000
$ws.Range("M75").Value = -1978
$ws.Range("N75").Value = -21872

$ws.Range("H78").Value = 11457
$ws.Range("I78").Value = 2914
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 8742
$ws.Range("L78").Value = 60000
$ws.Range("M78").Value = -4062
$ws.Range("N78").Value = -69360

$ws.Range("H80").Value = 590.0714
$ws.Range("I80").Value = 82
$ws.Range("J80").Value = 793.3
$ws.Range("K80").Value = 82
$ws.Range("L80").Value = 793.3
$ws.Range("M80").Value = 916
$ws.Range("N80").Value = -2789.3

$ws.Range("H83").Value = 590.0714
$ws.Range("I83").Value = 82
$ws.Range("J83").Value = 793.3
$ws.Range("K83").Value = 410
$ws.Range("L83").Value = 3966.5
$ws.Range("M83").Value = 4582
$ws.Range("N83").Value = -13950.5

$ws.Range("H99").Value = 1017.9
$ws.Range("I99").Value = 886.4706
$ws.Range("J99").Value = 1762.6666
$ws.Range("K99").Value = 886.4706
$ws.Range("L99").Value = 1762.6666
$ws.Range("M99").Value = 611.5294
$ws.Range("N99").Value = -4758.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 490.65518
$ws.Range("I22").Value = 474.125
$ws.Range("J22").Value = 570
$ws.Range("K22").Value = 474.125
$ws.Range("L22").Value = 570
$ws.Range("M22").Value = -124.125
$ws.Range("N22").Value = -1270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 218.33333
$ws.Range("I13").Value = 227.5
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 227.5
$ws.Range("L13").Value = 200
$ws.Range("M13").Value = -88.5
$ws.Range("N13").Value = -478

$ws.Range("H41").Value = 5533.3335
$ws.Range("I41").Value = 400
$ws.Range("J41").Value = 15800
$ws.Range("K41").Value = 400
$ws.Range("L41").Value = 15800
$ws.Range("M41").Value = -45
$ws.Range("N41").Value = -16510

$ws.Range("H132").Value = 4604.32
$ws.Range("I132").Value = 4413.875
$ws.Range("J132").Value = 4942.8887
$ws.Range("K132").Value = 13241.625
$ws.Range("L132").Value = 14828.6661
$ws.Range("M132").Value = -10711.625
$ws.Range("N132").Value = -19888.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 501750
$ws.Range("I136").Value = 335666.66
$ws.Range("J136").Value = 1000000
$ws.Range("K136").Value = 1006999.98
$ws.Range("L136").Value = 3000000
$ws.Range("M136").Value = -1004449.98
$ws.Range("N136").Value = -3005100
